# Commit: "rm print array and sheet; add test data in record.xlxs"
#
# Adds new test-data rows/cells to Sheet1 (the sheet's embedded
# printer-settings reference is dropped by the engine on save, which is
# the "rm print array" part of the commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Second date value next to the existing 2017-05-21 entry
$ws.Range("C2").Value = 20170522

# New test-data rows
$ws.Range("A4").Value = "胡蒙"
$ws.Range("C4").Value = "测试 "

$ws.Range("A5").Value = "胡蒙"
$ws.Range("C5").Value = "测试 "

$ws.Range("A6").Value = "郭闻浩"
$ws.Range("C6").Value = "测试 "
